# === Edit script: apply weight-tests changes to testResults.xlsx ===
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Cell($addr, $val) {
    $ws.Range($addr).Value2 = $val
}

# --- 1) Reformat existing commandRun strings (G column) to use space-separated,
#     quoted '-flag' args instead of the old 'flag=' dash-joined form.
#     Each commandRun string is shared by 6 rows, so every row in the block must
#     be rewritten identically to keep it as a single shared string.

# rows 109-114
foreach ($r in 109..114) {
    Set-Cell "G$r" 'testSamples24-7.py -batchSize=8 -epochs=100 -lr=0.001 -evalDetailLine=''no background tumor boxes adamw'' -hasBackground=f -usesLargestBox=f -segmentsMultiple=1 -dropoutRate=0.2 -grouped2D=f'
}

# rows 116-121
foreach ($r in 116..121) {
    Set-Cell "G$r" 'testSamples24-7.py -batchSize=8 -epochs=100 -lr=0.001 -evalDetailLine=''background global boxes adamw'' -hasBackground=t -usesLargestBox=t -segmentsMultiple=1 -dropoutRate=0.2 -grouped2D=f'
}

# rows 123-128
foreach ($r in 123..128) {
    Set-Cell "G$r" 'testSamples24-7.py -batchSize=8 -epochs=100 -lr=0.001 -evalDetailLine=''background tumor boxes adamw'' -hasBackground=t -usesLargestBox=f -segmentsMultiple=1 -dropoutRate=0.2 -grouped2D=f'
}

# rows 130-135
foreach ($r in 130..135) {
    Set-Cell "G$r" 'testSamples24-7.py -batchSize=8 -epochs=100 -lr=0.001 -evalDetailLine=''no background global boxes adamw'' -hasBackground=f -usesLargestBox=t -segmentsMultiple=1 -dropoutRate=0.2 -grouped2D=f'
}

# --- 2) Append the new "Weight tests" block: rows 138-154 ---

# Row 138
Set-Cell 'A138' 'Weight tests, viewing performance of legacy vs new weights'

# Row 139
Set-Cell 'A139' 'legacy weights'

# Row 140
Set-Cell 'A140' 'Tests/0--/foldn5'
Set-Cell 'B140' 100
Set-Cell 'C140' 8
Set-Cell 'D140' 0.001
Set-Cell 'E140' 0.2
Set-Cell 'F140' 0.01
Set-Cell 'G140' 'python testSamples26-7.py -batchSize=8 -epochs=100 -lr=0.001 -evalDetailLine=''legacy weights'' -hasBackground=f -usesLargestBox=f -segmentsMultiple=1 -dropoutRate=0.2 -grouped2D=f -modelChosen=''Small2DResnet'''
Set-Cell 'I140' '[None, None]'
Set-Cell 'L140' '[None, None]'
Set-Cell 'O140' '{0: 5.2, 1: 6.2, 2: 5.6}'
Set-Cell 'P140' 0.2705882352941177
Set-Cell 'Q140' 0.2772828478710832
Set-Cell 'R140' 0.2705882352941177
Set-Cell 'S140' 0.09799218018345172
Set-Cell 'T140' 0.1015540959098514
Set-Cell 'U140' 0.09799218018345172
Set-Cell 'V140' '[0.17647058823529413, 0.35294117647058826, 0.23529411764705882, 0.35294117647058826, 0.23529411764705882]'
Set-Cell 'W140' '[0.20168067226890754, 0.3539467068878834, 0.2184873949579832, 0.3770053475935829, 0.23529411764705882]'
Set-Cell 'X140' '[0.17647058823529413, 0.35294117647058826, 0.23529411764705882, 0.35294117647058826, 0.23529411764705882]'
Set-Cell 'Y140' '[14, 11, 13, 15, 15]'

# Row 141
Set-Cell 'A141' 'newer weights'

# Row 142
Set-Cell 'A142' 'Tests/0--/foldn5'
Set-Cell 'B142' 100
Set-Cell 'C142' 8
Set-Cell 'D142' 0.001
Set-Cell 'E142' 0.2
Set-Cell 'F142' 0.01
Set-Cell 'G142' 'python testSamples26-7.py -batchSize=8 -epochs=100 -lr=0.001 -evalDetailLine=''newer weights'' -hasBackground=f -usesLargestBox=f -segmentsMultiple=1 -dropoutRate=0.2 -grouped2D=f -modelChosen=''Small2DResnet'''
Set-Cell 'I142' '[None, None]'
Set-Cell 'L142' '[None, None]'
Set-Cell 'O142' '{0: 4.2, 1: 9.8, 2: 3.0}'
Set-Cell 'P142' 0.4941176470588236
Set-Cell 'Q142' 0.4576137175517981
Set-Cell 'R142' 0.4941176470588236
Set-Cell 'S142' 0.1326820391940614
Set-Cell 'T142' 0.1437866546233394
Set-Cell 'U142' 0.1326820391940614
Set-Cell 'V142' '[0.35294117647058826, 0.6470588235294118, 0.5294117647058824, 0.47058823529411764, 0.47058823529411764]'
Set-Cell 'W142' '[0.3288613691090472, 0.6039215686274512, 0.5409753645047762, 0.36463708290333674, 0.44967320261437904]'
Set-Cell 'X142' '[0.35294117647058826, 0.6470588235294118, 0.5294117647058824, 0.47058823529411764, 0.47058823529411764]'
Set-Cell 'Y142' '[19, 25, 35, 17, 21]'

# Row 143
Set-Cell 'A143' 'legacy weights'

# Row 144
Set-Cell 'A144' 'Tests/0--/foldn5'
Set-Cell 'B144' 100
Set-Cell 'C144' 8
Set-Cell 'D144' 0.001
Set-Cell 'E144' 0.2
Set-Cell 'F144' 0.01
Set-Cell 'G144' 'python testSamples26-7.py -batchSize=8 -epochs=100 -lr=0.001 -evalDetailLine=''legacy weights'' -hasBackground=f -usesLargestBox=f -segmentsMultiple=1 -dropoutRate=0.2 -grouped2D=f -modelChosen=''Large2DResnet'''
Set-Cell 'I144' '[None, None]'
Set-Cell 'L144' '[None, None]'
Set-Cell 'O144' '{0: 3.4, 1: 8.0, 2: 5.6}'
Set-Cell 'P144' 0.3411764705882353
Set-Cell 'Q144' 0.3331857719164221
Set-Cell 'R144' 0.3411764705882353
Set-Cell 'S144' 0.1890569159254115
Set-Cell 'T144' 0.1882956898952459
Set-Cell 'U144' 0.1890569159254115
Set-Cell 'V144' '[0.17647058823529413, 0.47058823529411764, 0.47058823529411764, 0.4117647058823529, 0.17647058823529413]'
Set-Cell 'W144' '[0.17647058823529413, 0.492436974789916, 0.47058823529411764, 0.34506051224317474, 0.18137254901960784]'
Set-Cell 'X144' '[0.17647058823529413, 0.47058823529411764, 0.47058823529411764, 0.4117647058823529, 0.17647058823529413]'
Set-Cell 'Y144' '[15, 15, 16, 28, 11]'

# Row 145
Set-Cell 'A145' 'newer weights'

# Row 146
Set-Cell 'A146' 'Tests/0--/foldn5'
Set-Cell 'B146' 100
Set-Cell 'C146' 8
Set-Cell 'D146' 0.001
Set-Cell 'E146' 0.2
Set-Cell 'F146' 0.01
Set-Cell 'G146' 'python testSamples26-7.py -batchSize=8 -epochs=100 -lr=0.001 -evalDetailLine=''newer weights'' -hasBackground=f -usesLargestBox=f -segmentsMultiple=1 -dropoutRate=0.2 -grouped2D=f -modelChosen=''Large2DResnet'''
Set-Cell 'I146' '[None, None]'
Set-Cell 'L146' '[None, None]'
Set-Cell 'O146' '{0: 2.6, 1: 8.4, 2: 6.0}'
Set-Cell 'P146' 0.3411764705882353
Set-Cell 'Q146' 0.3078145565961993
Set-Cell 'R146' 0.3411764705882353
Set-Cell 'S146' 0.1583449060198768
Set-Cell 'T146' 0.1523269651150921
Set-Cell 'U146' 0.1583449060198768
Set-Cell 'V146' '[0.29411764705882354, 0.35294117647058826, 0.35294117647058826, 0.17647058823529413, 0.5294117647058824]'
Set-Cell 'W146' '[0.273109243697479, 0.3025210084033613, 0.3740024681201152, 0.13003095975232196, 0.45940910300771887]'
Set-Cell 'X146' '[0.29411764705882354, 0.35294117647058826, 0.35294117647058826, 0.17647058823529413, 0.5294117647058824]'
Set-Cell 'Y146' '[17, 16, 18, 12, 22]'

# Row 148
Set-Cell 'A148' 'testing Padding on small model with tumor sizes (0 and valid padding seems to perform the same'

# Row 149
Set-Cell 'A149' 'newer weights with padding 0'

# Row 150
Set-Cell 'A150' 'Tests/0--/foldn5'
Set-Cell 'B150' 100
Set-Cell 'C150' 8
Set-Cell 'D150' 0.001
Set-Cell 'E150' 0.2
Set-Cell 'F150' 0.01
Set-Cell 'G150' 'python testSamples26-7.py -batchSize=8 -epochs=100 -lr=0.001 -evalDetailLine=''newer weights with padding 0'' -hasBackground=f -usesLargestBox=f -segmentsMultiple=1 -dropoutRate=0.2 -grouped2D=f -modelChosen=''Small2DResnet'''
Set-Cell 'I150' '[None, None]'
Set-Cell 'L150' '[None, None]'
Set-Cell 'O150' '{0: 4.2, 1: 9.8, 2: 3.0}'
Set-Cell 'P150' 0.4941176470588236
Set-Cell 'Q150' 0.4576137175517981
Set-Cell 'R150' 0.4941176470588236
Set-Cell 'S150' 0.1326820391940614
Set-Cell 'T150' 0.1437866546233394
Set-Cell 'U150' 0.1326820391940614
Set-Cell 'V150' '[0.35294117647058826, 0.6470588235294118, 0.5294117647058824, 0.47058823529411764, 0.47058823529411764]'
Set-Cell 'W150' '[0.3288613691090472, 0.6039215686274512, 0.5409753645047762, 0.36463708290333674, 0.44967320261437904]'
Set-Cell 'X150' '[0.35294117647058826, 0.6470588235294118, 0.5294117647058824, 0.47058823529411764, 0.47058823529411764]'
Set-Cell 'Y150' '[19, 25, 35, 17, 21]'

# Row 151
Set-Cell 'A151' 'newer weights with padding same'

# Row 152
Set-Cell 'A152' 'Tests/0--/foldn5'
Set-Cell 'B152' 100
Set-Cell 'C152' 8
Set-Cell 'D152' 0.001
Set-Cell 'E152' 0.2
Set-Cell 'F152' 0.01
Set-Cell 'G152' 'python testSamples26-7.py -batchSize=8 -epochs=100 -lr=0.001 -evalDetailLine=''newer weights with padding same'' -hasBackground=f -usesLargestBox=f -segmentsMultiple=1 -dropoutRate=0.2 -grouped2D=f -modelChosen=''Small2DResnet'''
Set-Cell 'I152' '[None, None]'
Set-Cell 'L152' '[None, None]'
Set-Cell 'O152' '{0: 2.4, 1: 8.4, 2: 6.2}'
Set-Cell 'P152' 0.3529411764705883
Set-Cell 'Q152' 0.329394780536649
Set-Cell 'R152' 0.3529411764705883
Set-Cell 'S152' 0.1032928274217746
Set-Cell 'T152' 0.08403464523761299
Set-Cell 'U152' 0.1032928274217746
Set-Cell 'V152' '[0.35294117647058826, 0.47058823529411764, 0.35294117647058826, 0.35294117647058826, 0.23529411764705882]'
Set-Cell 'W152' '[0.3585434173669468, 0.4125632153313814, 0.35173453996983406, 0.28431372549019607, 0.23981900452488686]'
Set-Cell 'X152' '[0.35294117647058826, 0.47058823529411764, 0.35294117647058826, 0.35294117647058826, 0.23529411764705882]'
Set-Cell 'Y152' '[28, 24, 19, 23, 16]'

# Row 153
Set-Cell 'A153' 'newer weights with padding valid'

# Row 154
Set-Cell 'A154' 'Tests/0--/foldn5'
Set-Cell 'B154' 100
Set-Cell 'C154' 8
Set-Cell 'D154' 0.001
Set-Cell 'E154' 0.2
Set-Cell 'F154' 0.01
Set-Cell 'G154' 'python testSamples26-7.py -batchSize=8 -epochs=100 -lr=0.001 -evalDetailLine=''newer weights with padding valid'' -hasBackground=f -usesLargestBox=f -segmentsMultiple=1 -dropoutRate=0.2 -grouped2D=f -modelChosen=''Small2DResnet'''
Set-Cell 'I154' '[None, None]'
Set-Cell 'L154' '[None, None]'
Set-Cell 'O154' '{0: 4.2, 1: 9.8, 2: 3.0}'
Set-Cell 'P154' 0.4941176470588236
Set-Cell 'Q154' 0.4576137175517981
Set-Cell 'R154' 0.4941176470588236
Set-Cell 'S154' 0.1326820391940614
Set-Cell 'T154' 0.1437866546233394
Set-Cell 'U154' 0.1326820391940614
Set-Cell 'V154' '[0.35294117647058826, 0.6470588235294118, 0.5294117647058824, 0.47058823529411764, 0.47058823529411764]'
Set-Cell 'W154' '[0.3288613691090472, 0.6039215686274512, 0.5409753645047762, 0.36463708290333674, 0.44967320261437904]'
Set-Cell 'X154' '[0.35294117647058826, 0.6470588235294118, 0.5294117647058824, 0.47058823529411764, 0.47058823529411764]'
Set-Cell 'Y154' '[19, 25, 35, 17, 21]'

# --- 3) Update sheet view state: scroll position + active selection ---
$ws.Range("A155").Select()
